$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shp = $s.Shapes.Item(3)
$tbl = $shp.Table
$tbl.ApplyStyle("{4DDFAC36-9B07-4B83-BC1A-3EF5311210D6}")
